# Fruta / hortaliza, semanal
# Insert a new data row at row 113 (shifting existing rows 113:191 down to
# 114:192) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 113; this shifts rows 113-191
# down to 114-192 and extends the used range to A1:T192.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new record.
$ws.Cells.Item(113, 1).Value  = 3
$ws.Cells.Item(113, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(113, 3).Value  = "Coquimbo"
$ws.Cells.Item(113, 4).Value  = 44596
$ws.Cells.Item(113, 5).Value  = 5
$ws.Cells.Item(113, 6).Value  = "Fruta"
$ws.Cells.Item(113, 7).Value  = 100101
$ws.Cells.Item(113, 8).Value  = "Berries"
$ws.Cells.Item(113, 9).Value  = 100101001
$ws.Cells.Item(113, 10).Value = "Arándano (blue)"
$ws.Cells.Item(113, 11).Value = "Sin especificar"
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 28
$ws.Cells.Item(113, 14).Value = 4000
$ws.Cells.Item(113, 15).Value = 4000
$ws.Cells.Item(113, 16).Value = 4000
$ws.Cells.Item(113, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(113, 18).Value = "Provincia de Linares"
$ws.Cells.Item(113, 19).Value = 2000
$ws.Cells.Item(113, 20).Value = 2
